$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- First summary table (rows 3-5, columns I:L) ---------------------------
# Header row copied from A1:D1 (same four labels / shared strings)
$ws.Range("I3").Value = "1 mg/ml"
$ws.Range("J3").Value = "0.1 mg/ml"
$ws.Range("K3").Value = "0.01 mg/ml"
$ws.Range("L3").Value = "Control"

# Row 4 - averages of each original data column (A2:A22, B2:B22, C2:C22, D2:D22)
$ws.Range("I4").Formula = "=AVERAGE(A2:A22)"
$ws.Range("J4").Formula = "=AVERAGE(B2:B22)"
$ws.Range("K4").Formula = "=AVERAGE(C2:C22)"
$ws.Range("L4").Formula = "=AVERAGE(D2:D22)"

# Row 5 - sample standard deviation of each original data column
$ws.Range("I5").Formula = "=STDEV.S(A2:A22)"
$ws.Range("J5").Formula = "=STDEV.S(B2:B22)"
$ws.Range("K5").Formula = "=STDEV.S(C2:C22)"
$ws.Range("L5").Formula = "=STDEV.S(D2:D22)"

# --- Second summary table (rows 8-10, columns I:L) -------------------------
# Same header repeated
$ws.Range("I8").Value = "1 mg/ml"
$ws.Range("J8").Value = "0.1 mg/ml"
$ws.Range("K8").Value = "0.01 mg/ml"
$ws.Range("L8").Value = "Control"

# Row 9 - averages normalized against the Control average ($L$4)
$ws.Range("I9").Formula = '=I4/$L$4'
$ws.Range("J9").Formula = '=J4/$L$4'
$ws.Range("K9").Formula = '=K4/$L$4'
$ws.Range("L9").Formula = '=L4/$L$4'

# Row 10 - standard deviations normalized against the Control average ($L$4)
$ws.Range("I10").Formula = '=I5/$L$4'
$ws.Range("J10").Formula = '=J5/$L$4'
$ws.Range("K10").Formula = '=K5/$L$4'
$ws.Range("L10").Formula = '=L5/$L$4'

# The normalized table is formatted as percentages
$ws.Range("I9:L10").Style = "Percent"

# Leave the selection on the first summary table, matching the saved view
$ws.Range("I3:L5").Select()
